# project-overview.pptx — "Add files via upload" re-export.
#
# The underlying change is: 8 straight-connector shapes (the arrow
# diagram overlay) were removed from the "Solution" slide (slide 6).
# Because the whole deck's shapes were minted from one sequential
# Google-Slides id counter, freeing those 8 ids (103-110) causes every
# later slide/notes-page shape (slides 7-12 and their notes pages) to
# be renumbered down by 8 on re-export. We reproduce the renumbering by
# rewriting each shape's Name (the id itself is assigned internally by
# PowerPoint and isn't settable through automation).

$p = $ppt.ActivePresentation

# 1) Remove the 8 connector ("arrow") shapes from slide 6.
$solutionSlide = $p.Slides.Item(6)
for ($i = $solutionSlide.Shapes.Count; $i -ge 1; $i--) {
    $shp = $solutionSlide.Shapes.Item($i)
    if ($shp.Name -like "Google Shape;10*;p18" -and $shp.Connector) {
        $shp.Delete()
    }
}

# 2) Shift every shape name on slides 7-12 (and their notes pages) down
#    by 8, keeping the "Google Shape;<id>;<suffix>" naming scheme intact.
function Rename-ShapesDown8($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $parts = $shp.Name -split ";"
        if ($parts.Count -eq 3 -and $parts[0] -eq "Google Shape") {
            $oldNum = [int]$parts[1]
            $newNum = $oldNum - 8
            $shp.Name = "Google Shape;" + $newNum + ";" + $parts[2]
        }
    }
}

for ($idx = 7; $idx -le 12; $idx++) {
    $slide = $p.Slides.Item($idx)
    Rename-ShapesDown8 $slide.Shapes
    Rename-ShapesDown8 $slide.NotesPage.Shapes
}
